$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 150 (shifts existing rows 150..235 down to 151..236,
# preserving row formatting/styles such as the date style on column D).
$ws.Rows("150:150").Insert()

# Populate the newly inserted row 150 with the new observation.
$ws.Cells.Item(150, 1).Value = 3
$ws.Cells.Item(150, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(150, 3).Value = "Coquimbo"
$ws.Cells.Item(150, 4).Value = 44606
$ws.Cells.Item(150, 5).Value = 5
$ws.Cells.Item(150, 6).Value = 100112001
$ws.Cells.Item(150, 7).Value = "Berenjena"
$ws.Cells.Item(150, 8).Value = "Sin especificar"
$ws.Cells.Item(150, 9).Value = "Primera"
$ws.Cells.Item(150, 10).Value = 95
$ws.Cells.Item(150, 11).Value = 10000
$ws.Cells.Item(150, 12).Value = 11000
$ws.Cells.Item(150, 13).Value = 10474
$ws.Cells.Item(150, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(150, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(150, 16).Value = 175
$ws.Cells.Item(150, 17).Value = 60
$ws.Cells.Item(150, 18).Value = "Hortaliza"
